$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.227.51"
$ws.Range("E2").Value = "  +6.39%  "
$ws.Range("D3").Value = "2.751.33"
$ws.Range("E3").Value = "  +5.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.87%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").Value = "2.784.33"
$ws.Range("E9").Value = "  +6.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("E11").Value = "  +8.84%  "
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.158"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "3.238.70"
$ws.Range("E14").Value = "  +5.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.15%  "
$ws.Range("D16").Value = "64.081.29"
$ws.Range("E16").Value = "  +6.16%  "
$ws.Range("E17").Value = "  +9.57%  "
$ws.Range("D18").Value = "2.776.58"
$ws.Range("E18").Value = "  +6.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.16%  "
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.537"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "0.0₃0921"
$ws.Range("E29").Value = "  +15.42%  "
$ws.Range("E30").Value = "  +6.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  +18.06%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.40%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +12.00%  "
$ws.Range("E37").Value = "  +11.35%  "
$ws.Range("E38").Value = "  +10.20%  "
$ws.Range("E39").Value = "  +21.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "350.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0595"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.26%  "
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0260"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.48%  "
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").Value = "2.184.47"
$ws.Range("E51").Value = "  +7.73%  "
